$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 597
$ws1.Cells.Item(7, 6).Value = 14921
$ws1.Cells.Item(9, 6).Value = 4
$ws1.Cells.Item(10, 6).Value = 675
$ws1.Cells.Item(11, 6).Value = 15223
$ws1.Cells.Item(12, 6).Value = 35
$ws1.Cells.Item(13, 6).Value = 8728
$ws1.Cells.Item(14, 6).Value = 333
$ws1.Cells.Item(17, 6).Value = 181
$ws1.Cells.Item(19, 6).Value = 0
$ws1.Cells.Item(20, 6).Value = 11
$ws1.Cells.Item(21, 6).Value = 17
$ws1.Cells.Item(23, 6).Value = 21
$ws1.Cells.Item(26, 6).Value = 1081
$ws1.Cells.Item(29, 6).Value = 56
$ws1.Cells.Item(33, 6).Value = 29
$ws1.Cells.Item(36, 6).Value = 270
$ws1.Cells.Item(39, 6).Value = 5361
$ws1.Cells.Item(40, 6).Value = 5225

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 597
$ws4.Cells.Item(7, 6).Value = 14921
$ws4.Cells.Item(9, 6).Value = 4
$ws4.Cells.Item(10, 6).Value = 675
$ws4.Cells.Item(11, 6).Value = 15223
$ws4.Cells.Item(12, 6).Value = 35
$ws4.Cells.Item(13, 6).Value = 8728
$ws4.Cells.Item(14, 6).Value = 333
$ws4.Cells.Item(18, 6).Value = 181
$ws4.Cells.Item(21, 6).Value = 11
$ws4.Cells.Item(22, 6).Value = 17
$ws4.Cells.Item(24, 6).Value = 21
$ws4.Cells.Item(27, 6).Value = 1081
$ws4.Cells.Item(30, 6).Value = 56
$ws4.Cells.Item(36, 6).Value = 29
$ws4.Cells.Item(39, 6).Value = 270
$ws4.Cells.Item(42, 6).Value = 5361
$ws4.Cells.Item(43, 6).Value = 5225
